$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.051.44'
$ws.Range("E2").Value = '  -3.63%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.443.66'
$ws.Range("E3").Value = '  -3.98%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '525.52'
$ws.Range("E5").Value = '  -2.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.00'
$ws.Range("E6").Value = '  -8.53%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.64%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.550'
$ws.Range("E8").Value = '  -4.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.443.86'
$ws.Range("E9").Value = '  -5.12%  '

$ws.Range("B10").Value = 'TRON'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.160'
$ws.Range("E10").Value = '  -0.39%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0976'
$ws.Range("E11").Value = '  -4.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.28'
$ws.Range("E12").Value = '  -3.72%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.339'
$ws.Range("E13").Value = '  -6.50%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.877.72'
$ws.Range("E14").Value = '  -3.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '57.937.72'
$ws.Range("E15").Value = '  -3.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.35'
$ws.Range("E16").Value = '  -8.29%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000137'
$ws.Range("E17").Value = '  -4.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.448.36'
$ws.Range("E18").Value = '  -5.36%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.56'
$ws.Range("E19").Value = '  -6.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '317.29'
$ws.Range("E20").Value = '  -3.39%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.14'
$ws.Range("E21").Value = '  -5.14%  '

$ws.Range("E22").Value = '  -0.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.64'
$ws.Range("E23").Value = '  -5.39%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.89'
$ws.Range("E24").Value = '  -1.96%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.403'
$ws.Range("E25").Value = '  -7.86%  '

$ws.Range("E26").Value = '  -3.01%  '

$ws.Range("E27").Value = '  -1.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.38'
$ws.Range("E28").Value = '  -8.20%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0742'
$ws.Range("E29").Value = '  -8.00%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.73'
$ws.Range("E30").Value = '  -4.51%  '

$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.41'
$ws.Range("E31").Value = '  -10.13%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '162.60'
$ws.Range("E32").Value = '  -1.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.18%  '

$ws.Range("E34").Value = '  -12.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.03'
$ws.Range("E35").Value = '  -4.19%  '

$ws.Range("E36").Value = '  -9.79%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.97'
$ws.Range("E37").Value = '  -11.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.52'
$ws.Range("E38").Value = '  -7.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.26'
$ws.Range("E39").Value = '  -2.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.48'
$ws.Range("E40").Value = '  -7.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.774'
$ws.Range("E41").Value = '  -7.83%  '

$ws.Range("E42").Value = '  +0.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '269.36'
$ws.Range("E43").Value = '  -11.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.95'
$ws.Range("E44").Value = '  -11.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.85'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.582'
$ws.Range("E46").Value = '  -4.89%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0914'
$ws.Range("E47").Value = '  -2.72%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.25'
$ws.Range("E48").Value = '  -6.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0499'
$ws.Range("E49").Value = '  -4.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0216'
$ws.Range("E50").Value = '  -6.55%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.74'
$ws.Range("E51").Value = '  -8.82%  '
